$wb = $excel.ActiveWorkbook

$wsCover   = $wb.Worksheets.Item("Cover")
$wsRegions = $wb.Worksheets.Item("Regions")
$wsPWR     = $wb.Worksheets.Item("PWR")
$wsTRA     = $wb.Worksheets.Item("TRA")
$wsSRV     = $wb.Worksheets.Item("SRV")

# ---------------------------------------------------------------------------
# Data / assumption changes
# ---------------------------------------------------------------------------

# TRA: CAR/LGT max-growth share input used to derive L6:L8 (-$B$53*L27/1000)
$wsTRA.Range("L27").Value = 0.225
$wsTRA.Range("L28").Value = 0.225

# TRA: annual growth-rate assumption feeding the projection table (C50:O53)
$wsTRA.Range("C47").Value = 0.15

# SRV: Biomass / Biogas max-growth-rate starting values (feed L11/L12 = -D17/-D18)
$wsSRV.Range("D17").Value = 0.1
$wsSRV.Range("D18").Value = 0.1

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------

# Update the remembered selection on each sheet that changed, then leave PWR
# active last so it ends up the selected tab (matches activeTab in the diff).
$wsTRA.Activate()
$wsTRA.Range("L29").Select() | Out-Null
$winTRA = $excel.ActiveWindow
$winTRA.Zoom = 85

$wsSRV.Activate()
$wsSRV.Range("D18").Select() | Out-Null

$wsPWR.Activate()
$wsPWR.Range("K22").Select() | Out-Null
